$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Best-effort: VBA code names (cosmetic; Excel assigns these on every
# save). Harmless if the host doesn't persist them. ---
try { $wb.CodeName = "EstaPastaDeTrabalho" } catch { }
try { $ws.CodeName = "Planilha1" } catch { }

# --- Insert a new column between "Classe" (I) and "Conta Contábil" (old J)
# for the new "Desc. Classe" field, shifting everything after it right by one. ---
$ws.Columns("J:J").Insert()

# --- Normalize / retitle the header row (accents and casing cleaned up,
# plus the brand-new "Desc. Classe" column). ---
$ws.Range("A1").Value = "Nº Imobilizado"
$ws.Range("B1").Value = "Sub Nº"
$ws.Range("C1").Value = "Data Inicio da Depreciacao"
$ws.Range("D1").Value = "Descricao"
$ws.Range("E1").Value = "Valor Aquisicao"
$ws.Range("F1").Value = "Depreciacao Acum."
$ws.Range("G1").Value = "Valor Contabil"
$ws.Range("H1").Value = "Centro Custos"
$ws.Range("I1").Value = "Classe"
$ws.Range("J1").Value = "Desc. Classe"
$ws.Range("K1").Value = "Conta Contabil"
$ws.Range("L1").Value = "Desc. Conta Contabil"
$ws.Range("M1").Value = "Vida Util (anos)"
$ws.Range("N1").Value = "Vida Util (períodos)"
$ws.Range("O1").Value = "Data Fim Depreciacao"
$ws.Range("P1").Value = "Auxiliar 1"
$ws.Range("Q1").Value = "Auxiliar 2"

# --- Column width for the newly inserted column (raw OOXML width 11,
# i.e. ColumnWidth = 11 - 5/6 once Excel's fixed character-width offset
# is applied). ---
$ws.Columns("J:J").ColumnWidth = 10.16666666666667

# --- View/selection: the sheet is now scrolled/selected further left,
# landing on F8. ---
$ws.Range("F8").Select()
